$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "expected" value text for the existing Facebook rows.
$ws.Range("D2").Value = "Forgotten Password | Can't Log In | Facebook"
$ws.Range("D3").Value = "Forgotten Password | Can't Log In | Facebook"
$ws.Range("D4").Value = "Forgotten Password | Can't Log In | Facebook"

# Add the new "paymentTest" row.
$ws.Range("A5").Value = "paymentTest"
$ws.Range("B5").Value = "Ron"
$ws.Range("C5").Value = "Ron123"
$ws.Range("D5").Value = "Forgotten Password | Can't Log In | Facebook"

# Match the styling used by the other test-case name cells (A2/A3).
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
